$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: new log entry ("Implemented feedback" / 1.5 hours / "Implemented some missed ui")
$ws.Range("A13").Value = "Implemented feedback"
$ws.Range("B13").Value = 1.5
$ws.Range("D13").Value = "Implemented some missed ui"

# Move the saved selection to match the author's final cursor position
$ws.Range("L15").Select()
